$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 978.1429000000001
$ws.Range("I8").Value = 724.5
$ws.Range("K8").Value = 2173.5
$ws.Range("M8").Value = -2034.5

$ws.Range("H9").Value = 300
$ws.Range("I9").Value = 300
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 300
$ws.Range("L9").Value = 300
$ws.Range("M9").Value = -131
$ws.Range("N9").Value = -638

$ws.Range("H10").Value = 70005
$ws.Range("J10").Value = 70005
$ws.Range("L10").Value = 70005
$ws.Range("N10").Value = -70591

$ws.Range("H19").Value = 780.65216
$ws.Range("I19").Value = 753.375
$ws.Range("J19").Value = 795.2
$ws.Range("K19").Value = 753.375
$ws.Range("L19").Value = 795.2
$ws.Range("M19").Value = -578.375
$ws.Range("N19").Value = -1145.2

$ws.Range("H48").Value = 125702.125
$ws.Range("I48").Value = 333739
$ws.Range("J48").Value = 880
$ws.Range("K48").Value = 1001217
$ws.Range("L48").Value = 2640
$ws.Range("M48").Value = -1000925
$ws.Range("N48").Value = -3224

$ws.Range("H56").Value = 125702.125
$ws.Range("I56").Value = 333739
$ws.Range("J56").Value = 880
$ws.Range("K56").Value = 1001217
$ws.Range("L56").Value = 2640
$ws.Range("M56").Value = -1000683
$ws.Range("N56").Value = -3708

$ws.Range("H86").Value = 81261.2
$ws.Range("I86").Value = 1668.6666
$ws.Range("J86").Value = 200650
$ws.Range("K86").Value = 1668.6666
$ws.Range("L86").Value = 200650
$ws.Range("M86").Value = -545.6666
$ws.Range("N86").Value = -202896

$ws.Range("H89").Value = 81261.2
$ws.Range("I89").Value = 1668.6666
$ws.Range("J89").Value = 200650
$ws.Range("K89").Value = 8343.333000000001
$ws.Range("L89").Value = 1003250
$ws.Range("M89").Value = -2727.333000000001
$ws.Range("N89").Value = -1014482

$ws.Range("H138").Value = 2297.3403
$ws.Range("I138").Value = 1528.4445
$ws.Range("J138").Value = 3335.35
$ws.Range("K138").Value = 4585.333500000001
$ws.Range("L138").Value = 10006.05
$ws.Range("M138").Value = 554.6664999999994
$ws.Range("N138").Value = -20286.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2000
$ws.Range("I3").Value = 2000
$ws.Range("K3").Value = 2000
$ws.Range("M3").Value = -1885

$ws.Range("H61").Value = 1050.8572
$ws.Range("I61").Value = 649.6667
$ws.Range("J61").Value = 1351.75
$ws.Range("K61").Value = 649.6667
$ws.Range("L61").Value = 1351.75
$ws.Range("M61").Value = -437.6667
$ws.Range("N61").Value = -1775.75

$ws.Range("H136").Value = 1050.8572
$ws.Range("I136").Value = 649.6667
$ws.Range("J136").Value = 1351.75
$ws.Range("K136").Value = 1949.0001
$ws.Range("L136").Value = 4055.25
$ws.Range("M136").Value = 600.9999
$ws.Range("N136").Value = -9155.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 200.5
$ws.Range("I8").Value = 301
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 301
$ws.Range("L8").Value = 100
$ws.Range("M8").Value = -161
$ws.Range("N8").Value = -380

$ws.Range("H134").Value = 1712.579
$ws.Range("I134").Value = 1485.8125
$ws.Range("J134").Value = 2922
$ws.Range("K134").Value = 4457.4375
$ws.Range("L134").Value = 8766
$ws.Range("M134").Value = -1922.4375
$ws.Range("N134").Value = -13836

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6428.129
$ws.Range("I31").Value = 3004.875
$ws.Range("J31").Value = 10079.6
$ws.Range("K31").Value = 3004.875
$ws.Range("L31").Value = 10079.6
$ws.Range("M31").Value = -2709.875
$ws.Range("N31").Value = -10669.6

$ws.Range("H34").Value = 6428.129
$ws.Range("I34").Value = 3004.875
$ws.Range("J34").Value = 10079.6
$ws.Range("K34").Value = 3004.875
$ws.Range("L34").Value = 10079.6
$ws.Range("M34").Value = -2802.875
$ws.Range("N34").Value = -10483.6

$ws.Range("H107").Value = 694.0952
$ws.Range("I107").Value = 291.57144
$ws.Range("J107").Value = 895.3570999999999
$ws.Range("K107").Value = 291.57144
$ws.Range("L107").Value = 895.3570999999999
$ws.Range("M107").Value = 1628.42856
$ws.Range("N107").Value = -4735.3571

$ws.Range("H132").Value = 1693.3877
$ws.Range("I132").Value = 767.75757
$ws.Range("J132").Value = 3602.5
$ws.Range("K132").Value = 2303.27271
$ws.Range("L132").Value = 10807.5
$ws.Range("M132").Value = 226.7272899999998
$ws.Range("N132").Value = -15867.5

$ws.Range("H138").Value = 33656.555
$ws.Range("J138").Value = 33656.555
$ws.Range("L138").Value = 33656.555
$ws.Range("N138").Value = -43936.555

$ws.Range("H140").Value = 41407.2
$ws.Range("J140").Value = 41407.2
$ws.Range("L140").Value = 41407.2
$ws.Range("N140").Value = -51767.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 140.90909
$ws.Range("J2").Value = 147.88889
$ws.Range("L2").Value = 147.88889
$ws.Range("N2").Value = -373.88889

$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()

$ws.Range("H7").Value = 4381250
$ws.Range("J7").Value = 50000
$ws.Range("L7").Value = 50000
$ws.Range("N7").Value = -50224

$ws.Range("H8").Value = 4381250
$ws.Range("J8").Value = 50000
$ws.Range("L8").Value = 50000
$ws.Range("N8").Value = -50278

$ws.Range("H10").Value = 3334666.8
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8447.111000000001
$ws.Range("I132").Value = 9187.483
$ws.Range("J132").Value = 6807.7144
$ws.Range("K132").Value = 27562.449
$ws.Range("L132").Value = 20423.1432
$ws.Range("M132").Value = -25032.449
$ws.Range("N132").Value = -25483.1432

$ws.Range("H136").Value = 9806929
$ws.Range("I136").Value = 3135.7083
$ws.Range("J136").Value = 33336034
$ws.Range("K136").Value = 9407.124899999999
$ws.Range("L136").Value = 100008102
$ws.Range("M136").Value = -6857.124899999999
$ws.Range("N136").Value = -100013202

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 200
$ws.Range("I8").Value = 200
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 200
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -60
$ws.Range("N8").ClearContents()

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("I96").Value = 1500
$ws.Range("K96").Value = 1500
$ws.Range("M96").Value = -127

$ws.Range("H132").Value = 2062.7437
$ws.Range("I132").Value = 1901.3
$ws.Range("J132").Value = 2600.889
$ws.Range("K132").Value = 5703.9
$ws.Range("L132").Value = 7802.667
$ws.Range("M132").Value = -3173.9
$ws.Range("N132").Value = -12862.667

$ws.Range("H136").Value = 1951.0588
$ws.Range("I136").Value = 651.0333000000001
$ws.Range("J136").Value = 11701.25
$ws.Range("K136").Value = 1953.0999
$ws.Range("L136").Value = 35103.75
$ws.Range("M136").Value = 596.9000999999998
$ws.Range("N136").Value = -40203.75
